$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Last status check on: 24.02.2022 12:15"

$ws.Range("B9").Value = 38.9
$ws.Range("C9").Value = 38.5

$ws.Range("D9").Value = "'+0.4"
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = "2022-02-24 12:15:24"
